$wb = $excel.ActiveWorkbook

# --- Sheet 1: "phylip-programs" ---
$ws1 = $wb.Worksheets.Item("phylip-programs")

# Row 34 (seqboot.exe): add Rseqboot function name + date done
$ws1.Range("B34").Value = "Rseqboot"
$ws1.Range("C34").Value = 41723
$ws1.Range("C32").Copy()
$ws1.Range("C34").PasteSpecial(-4122)

# --- Sheet 2: "addt'l-functions" ---
$ws2 = $wb.Worksheets.Item("addt'l-functions")

$ws2.Range("A21").Value = "read.multi.dna"
$ws2.Range("B21").Value = 41723
$ws2.Range("C21").Value = "Yes"
$ws2.Range("D21").Value = "Reads multiple DNA datasets from file."

$ws2.Range("A22").Value = "read.multi.phylip.data"
$ws2.Range("B22").Value = 41723
$ws2.Range("C22").Value = "Yes"
$ws2.Range("D22").Value = "Reads multiple phylip.data datasets from file."

$ws2.Range("A23").Value = "read.multi.rest.data"
$ws2.Range("B23").Value = 41723
$ws2.Range("C23").Value = "Yes"
$ws2.Range("D23").Value = "Reads multiple rest.data datasets from file."

$ws2.Range("A24").Value = "read.phylip.data"
$ws2.Range("B24").Value = 41723
$ws2.Range("C24").Value = "Yes"
$ws2.Range("D24").Value = "Reads phylip.data."

$ws2.Range("A25").Value = "read.rest.data"
$ws2.Range("B25").Value = 41723
$ws2.Range("C25").Value = "Yes"
$ws2.Range("D25").Value = "Reads rest.data."

# copy date format from an existing date cell in same column
$ws2.Range("B20").Copy()
$ws2.Range("B21:B25").PasteSpecial(-4122)

# Column A width grew to fit the new, longer function names
$ws2.Columns.Item(1).ColumnWidth = 20.3

# Restore the selection on sheet 1 (last action performed == last active sheet/selection)
$ws1.Range("B35").Select()
